$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: C20 = " " (single space) -- set first so this shared string is registered
# right after "deleteCustomer", matching the original author's save order.
$ws.Range("C20").Value = " "

# Row 3: C3 "none" -> "This is the first test"
$ws.Range("C3").Value = "This is the first test"
# Row 4: C4 "none" -> "This is the second test"
$ws.Range("C4").Value = "This is the second test"
# Row 5: C5 "this is the first test" -> "This is the third test"
$ws.Range("C5").Value = "This is the third test"

# Row 13: A13 "cus_IBVovXeEyeRqgn" -> "cus_IBc0ERhRyxXWsL"
$ws.Range("A13").Value = "cus_IBc0ERhRyxXWsL"
# Rows 14-18: new customer ids
$ws.Range("A14").Value = "cus_IBc00mEJZdW8Kg"
$ws.Range("A15").Value = "cus_IBc09g64O3FaQE"
$ws.Range("A16").Value = "cus_IBc0ej42CA1Txb"
$ws.Range("A17").Value = "cus_IBc0J46XloVal7"
$ws.Range("A18").Value = "cus_IBc0swmV0KXVgB"

# Update selection to H15 (as seen in diff)
$ws.Range("H15").Select()
